# Gantt chart update: fill in "Actual Start" / "Actual Duration" values for
# the remaining tasks and mark them as 100% complete (G column = 1), clear
# the stray formatted-but-empty H5 cell, and update the current selection /
# scroll position to match where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Row 22: 3.1.1 Data cleaning & preprocessing ---
$ws.Range("E22").Value = 37
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 1

# --- Row 23: 3.1.2 Data sorting & info extraction ---
$ws.Range("E23").Value = 40
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 1

# --- Row 24: 3.1.3 Exploratory Data Analysis ---
$ws.Range("E24").Value = 41
$ws.Range("F24").Value = 6
$ws.Range("G24").Value = 1

# --- Row 26: 3.2.1 Create interface ---
$ws.Range("E26").Value = 36
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 1

# --- Row 27: 3.2.2 Database integration ---
$ws.Range("E27").Value = 43
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 1

# --- Row 28: 3.2.3 Data Insight Visualization ---
$ws.Range("E28").Value = 46
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 1

# --- Row 29: 3.3 Requirements Acceptance Testing ---
$ws.Range("E29").Value = 51
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 1

# --- Row 31: 4.1 Status reports ---
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 56
$ws.Range("G31").Value = 1

# --- Row 32: 4.2 Update plans and schedule ---
$ws.Range("E32").Value = 45
$ws.Range("F32").Value = 16
$ws.Range("G32").Value = 1

# --- Row 34: 5.1 Lesson learned ---
$ws.Range("E34").Value = 55
$ws.Range("F34").Value = 3
$ws.Range("G34").Value = 1

# --- Row 35: 5.2 Prepare final report ---
$ws.Range("E35").Value = 55
$ws.Range("F35").Value = 6
$ws.Range("G35").Value = 1

# Stray empty-but-styled cell left over in row 5 - remove it entirely.
$ws.Range("H5").Clear()

# Leave the view where the author last left it: scrolled down so row 19 is
# at the top, with F25 selected.
$ws.Range("F25").Select()
